# Append the new data row (id=6) to the "Daily APR" log on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "2025-09-02T09:30"
$ws.Range("C7").Value = 1.6718133009139704
